# Atualiza notas dos alunos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the previously-empty "C4" column (I) grades for each student
$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("I4").Value = 2
$ws.Range("I5").Value = 0
$ws.Range("I6").Value = 0

# Record the "Conceito" (RF) for the first student in column N
$ws.Range("N2").Value = "RF"

# Update the selection to match the saved workbook state
$ws.Range("N3").Select()
